$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp shown in the title cell ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 09:52"

# --- Update per-country COVID statistics with freshly reported figures ---
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes
$countryRange = $ws.Range("A4:A216")

function Set-CountryRow {
    param($CountryName, $B, $C, $D, $E, $F, $G, $H)

    $cell = $countryRange.Find($CountryName)
    $r = $cell.Row
    $ws.Cells.Item($r, 2).Value = $B
    $ws.Cells.Item($r, 3).Value = $C
    $ws.Cells.Item($r, 4).Value = $D
    $ws.Cells.Item($r, 5).Value = $E
    $ws.Cells.Item($r, 6).Value = $F
    $ws.Cells.Item($r, 7).Value = $G
    $ws.Cells.Item($r, 8).Value = $H
}

Set-CountryRow "Rusia"     57999 5236 4420 53066 700 57 513
Set-CountryRow "Singapur"  10141 1016 839  9291  23  0  11
Set-CountryRow "Moldavia"  2614  0    560  1981  212 1  73
Set-CountryRow "Armenia"   1473  72   633  816   30  0  24
Set-CountryRow "Georgia"   411   3    98   308   6   1  5
Set-CountryRow "Rumania"   9242  0    2153 6582  245 9  507
Set-CountryRow "Sri Lanka" 310   0    104  199   2   0  7

# --- Re-sort the country table (rows 4-216) descending by "Casos totales" (col B) ---
# so the ranking stays consistent after the above updates, exactly like the source
# feed re-exports the whole table sorted by total cases each refresh.
$sortRange = $ws.Range("A4:H216")
$sortKey = $ws.Range("B4:B216")
$sortRange.Sort($sortKey, 2)
